$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 2275.7273
$ws.Range("I4").Value = 1379.25
$ws.Range("J4").Value = 4666.3335
$ws.Range("K4").Value = 1379.25
$ws.Range("L4").Value = 4666.3335
$ws.Range("M4").Value = -1265.25
$ws.Range("N4").Value = -4894.3335
# Row 33
$ws.Range("H33").Value = 252.22223
$ws.Range("I33").Value = 195.66667
$ws.Range("K33").Value = 195.66667
$ws.Range("M33").Value = 33.33332999999999
# Row 51
$ws.Range("H51").Value = 5588.8887
$ws.Range("I51").Value = 6966.6665
$ws.Range("J51").Value = 4900
$ws.Range("K51").Value = 6966.6665
$ws.Range("L51").Value = 4900
$ws.Range("M51").Value = -6482.6665
$ws.Range("N51").Value = -5868
# Row 53
$ws.Range("H53").Value = 283.55
$ws.Range("I53").Value = 224.25
$ws.Range("J53").Value = 372.5
$ws.Range("K53").Value = 224.25
$ws.Range("L53").Value = 372.5
$ws.Range("M53").Value = 412.75
$ws.Range("N53").Value = -1646.5
# Row 62
$ws.Range("H62").Value = 7967.375
$ws.Range("J62").Value = 7967.375
$ws.Range("L62").Value = 7967.375
$ws.Range("N62").Value = -9215.375
# Row 65
$ws.Range("H65").Value = 7967.375
$ws.Range("J65").Value = 7967.375
$ws.Range("L65").Value = 39836.875
$ws.Range("N65").Value = -46076.875
# Row 70
$ws.Range("H70").Value = 6899.6
$ws.Range("I70").Value = 1100
$ws.Range("J70").Value = 8349.5
$ws.Range("K70").Value = 3300
$ws.Range("L70").Value = 25048.5
$ws.Range("M70").Value = -3030
$ws.Range("N70").Value = -25588.5
# Row 73
$ws.Range("H73").Value = 6899.6
$ws.Range("I73").Value = 1100
$ws.Range("J73").Value = 8349.5
$ws.Range("K73").Value = 3300
$ws.Range("L73").Value = 25048.5
$ws.Range("M73").Value = -2364
$ws.Range("N73").Value = -26920.5
# Row 82
$ws.Range("H82").Value = 999
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
# Row 85
$ws.Range("H85").Value = 999
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
# Row 96
$ws.Range("H96").Value = 5757.4736
$ws.Range("I96").Value = 7720.2856
$ws.Range("K96").Value = 23160.8568
$ws.Range("M96").Value = -21787.8568
# Row 98
$ws.Range("H98").Value = 3722.111
$ws.Range("I98").Value = 1002.5
$ws.Range("J98").Value = 5897.8
$ws.Range("K98").Value = 1002.5
$ws.Range("L98").Value = 5897.8
$ws.Range("M98").Value = 495.5
$ws.Range("N98").Value = -8893.799999999999
# Row 122
$ws.Range("H122").Value = 3722.111
$ws.Range("I122").Value = 1002.5
$ws.Range("J122").Value = 5897.8
$ws.Range("K122").Value = 3007.5
$ws.Range("L122").Value = 17693.4
$ws.Range("M122").Value = -557.5
$ws.Range("N122").Value = -22593.4
# Row 141
$ws.Range("H141").Value = 4094
$ws.Range("I141").Value = 4094
$ws.Range("K141").Value = 12282
$ws.Range("M141").Value = -7102

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4738.35
$ws.Range("I61").Value = 1727.5
$ws.Range("K61").Value = 1727.5
$ws.Range("M61").Value = -1515.5
# Row 63
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314
# Row 66
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568
# Row 82
$ws.Range("H82").Value = 44949.5
$ws.Range("J82").Value = 44949.5
$ws.Range("L82").Value = 44949.5
$ws.Range("N82").Value = -45671.5
# Row 85
$ws.Range("H85").Value = 44949.5
$ws.Range("J85").Value = 44949.5
$ws.Range("L85").Value = 44949.5
$ws.Range("N85").Value = -47445.5
# Row 132
$ws.Range("H132").Value = 1739.8276
$ws.Range("I132").Value = 1709.2963
$ws.Range("K132").Value = 5127.8889
$ws.Range("M132").Value = -2597.8889
# Row 136
$ws.Range("H136").Value = 4738.35
$ws.Range("I136").Value = 1727.5
$ws.Range("K136").Value = 5182.5
$ws.Range("M136").Value = -2632.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 766.875
$ws.Range("I86").Value = 766.875
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 766.875
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 356.125
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 766.875
$ws.Range("I89").Value = 766.875
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 3834.375
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 1781.625
$ws.Range("N89").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2416.3635
$ws.Range("I16").Value = 2183.1667
$ws.Range("K16").Value = 2183.1667
$ws.Range("M16").Value = -1896.1667
# Row 22
$ws.Range("H22").Value = 36286
$ws.Range("J22").Value = 62938.5
$ws.Range("L22").Value = 62938.5
$ws.Range("N22").Value = -63638.5
# Row 31
$ws.Range("H31").Value = 1685.0385
$ws.Range("I31").Value = 1455.25
$ws.Range("K31").Value = 1455.25
$ws.Range("M31").Value = -1160.25
# Row 34
$ws.Range("H34").Value = 1685.0385
$ws.Range("I34").Value = 1455.25
$ws.Range("K34").Value = 1455.25
$ws.Range("M34").Value = -1253.25
# Row 113
$ws.Range("H113").Value = 2416.3635
$ws.Range("I113").Value = 2183.1667
$ws.Range("K113").Value = 2183.1667
$ws.Range("M113").Value = -13.16670000000022
# Row 132
$ws.Range("H132").Value = 1535.5758
$ws.Range("I132").Value = 1232.5769
$ws.Range("K132").Value = 3697.7307
$ws.Range("M132").Value = -1167.7307
# Row 133
$ws.Range("H133").Value = 46648.75
$ws.Range("J133").Value = 46648.75
$ws.Range("L133").Value = 46648.75
$ws.Range("N133").Value = -51708.75
# Row 134
$ws.Range("H134").Value = 2477.111
$ws.Range("I134").Value = 2307
$ws.Range("K134").Value = 6921
$ws.Range("M134").Value = -4386

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 221200
$ws.Range("I2").Value = 183433.5
$ws.Range("J2").Value = 277849.75
$ws.Range("K2").Value = 1100601
$ws.Range("L2").Value = 1667098.5
$ws.Range("M2").Value = -1100488
$ws.Range("N2").Value = -1667324.5
# Row 92
$ws.Range("H92").Value = 251.69698
$ws.Range("J92").Value = 465.22223
$ws.Range("L92").Value = 1395.66669
$ws.Range("N92").Value = -3891.66669
# Row 94
$ws.Range("H94").Value = 15571.286
$ws.Range("J94").Value = 19800
$ws.Range("L94").Value = 59400
$ws.Range("N94").Value = -60752
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 94
$ws.Range("H94").Value = 39999
$ws.Range("J94").Value = 39999
$ws.Range("L94").Value = 39999
$ws.Range("N94").Value = -41351
# Row 113
$ws.Range("H113").Value = 1589.6666
$ws.Range("I113").Value = 1538.375
$ws.Range("K113").Value = 1538.375
$ws.Range("M113").Value = 631.625

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1450
$ws.Range("I22").Value = 1400
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -1105
$ws.Range("N22").Value = -2090
# Row 27
$ws.Range("H27").Value = 1450
$ws.Range("I27").Value = 1400
$ws.Range("J27").Value = 1500
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 1500
$ws.Range("M27").Value = -1293
$ws.Range("N27").Value = -1714
# Row 32
$ws.Range("H32").Value = 999
$ws.Range("I32").Value = 999
$ws.Range("K32").Value = 999
$ws.Range("M32").Value = -682
# Row 46
$ws.Range("H46").Value = 168216.67
# Row 61
$ws.Range("H61").Value = 4998
$ws.Range("I61").Value = 5137.8
$ws.Range("J61").Value = 4648.5
$ws.Range("K61").Value = 5137.8
$ws.Range("L61").Value = 4648.5
$ws.Range("M61").Value = -4935.8
$ws.Range("N61").Value = -5052.5
# Row 74
$ws.Range("H74").Value = 39548.5
$ws.Range("I74").Value = 39548.5
$ws.Range("K74").Value = 39548.5
$ws.Range("M74").Value = -38550.5
# Row 77
$ws.Range("H77").Value = 39548.5
$ws.Range("I77").Value = 39548.5
$ws.Range("K77").Value = 118645.5
$ws.Range("M77").Value = -113653.5
# Row 113
$ws.Range("H113").Value = 4998
$ws.Range("I113").Value = 5137.8
$ws.Range("J113").Value = 4648.5
$ws.Range("K113").Value = 5137.8
$ws.Range("L113").Value = 4648.5
$ws.Range("M113").Value = -2967.8
$ws.Range("N113").Value = -8988.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 75
$ws.Range("H75").Value = 87118
$ws.Range("I75").Value = 87118
$ws.Range("K75").Value = 87118
$ws.Range("M75").Value = -86182
# Row 78
$ws.Range("H78").Value = 87118
$ws.Range("I78").Value = 87118
$ws.Range("K78").Value = 261354
$ws.Range("M78").Value = -256674
# Row 132
$ws.Range("H132").Value = 1212.4667
$ws.Range("J132").Value = 993
$ws.Range("L132").Value = 2979
$ws.Range("N132").Value = -8039
